# Auto-generated: apply numeric cell updates per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 342.33334
$ws.Range("I18").Value = 384
$ws.Range("J18").Value = 259
$ws.Range("K18").Value = 384
$ws.Range("L18").Value = 259
$ws.Range("M18").Value = -100
$ws.Range("N18").Value = -827
# Row 100
$ws.Range("H100").Value = 5963.357
$ws.Range("I100").Value = 4813.4287
$ws.Range("J100").Value = 7113.2856
$ws.Range("K100").Value = 4813.4287
$ws.Range("L100").Value = 7113.2856
$ws.Range("M100").Value = -4272.4287
$ws.Range("N100").Value = -8195.285599999999
# Row 106
$ws.Range("H106").Value = 6040.3125
$ws.Range("I106").Value = 5831.7856
$ws.Range("K106").Value = 5831.7856
$ws.Range("M106").Value = -5200.7856
# Row 132
$ws.Range("H132").Value = 17803.588
$ws.Range("I132").Value = 2618.6785
$ws.Range("K132").Value = 7856.0355
$ws.Range("M132").Value = -5326.0355
# Row 140
$ws.Range("H140").Value = 39714.285
$ws.Range("J140").Value = 39714.285
$ws.Range("L140").Value = 39714.285
$ws.Range("N140").Value = -50074.285

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4393.4
$ws.Range("I61").Value = 4787.864
$ws.Range("J61").Value = 3308.625
$ws.Range("K61").Value = 4787.864
$ws.Range("L61").Value = 3308.625
$ws.Range("M61").Value = -4575.864
$ws.Range("N61").Value = -3732.625
# Row 97
$ws.Range("H97").Value = 2551.1765
$ws.Range("I97").Value = 1186.3572
$ws.Range("K97").Value = 1186.3572
$ws.Range("M97").Value = -690.3571999999999
# Row 132
$ws.Range("H132").Value = 3256.348
$ws.Range("I132").Value = 3134.6667
$ws.Range("J132").Value = 3694.4
$ws.Range("K132").Value = 9404.000100000001
$ws.Range("L132").Value = 11083.2
$ws.Range("M132").Value = -6874.000100000001
$ws.Range("N132").Value = -16143.2
# Row 136
$ws.Range("H136").Value = 4393.4
$ws.Range("I136").Value = 4787.864
$ws.Range("J136").Value = 3308.625
$ws.Range("K136").Value = 14363.592
$ws.Range("L136").Value = 9925.875
$ws.Range("M136").Value = -11813.592
$ws.Range("N136").Value = -15025.875

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 11906311
$ws.Range("I86").Value = 13159002
$ws.Range("J86").Value = 5748.5
$ws.Range("K86").Value = 13159002
$ws.Range("L86").Value = 5748.5
$ws.Range("M86").Value = -13157879
$ws.Range("N86").Value = -7994.5
# Row 89
$ws.Range("H89").Value = 11906311
$ws.Range("I89").Value = 13159002
$ws.Range("J89").Value = 5748.5
$ws.Range("K89").Value = 65795010
$ws.Range("L89").Value = 28742.5
$ws.Range("M89").Value = -65789394
$ws.Range("N89").Value = -39974.5
# Row 105
$ws.Range("H105").Value = 1872.7858
$ws.Range("I105").Value = 749.4761999999999
$ws.Range("K105").Value = 749.4761999999999
$ws.Range("M105").Value = 997.5238000000001
# Row 134
$ws.Range("H134").Value = 1879.2122
$ws.Range("I134").Value = 1293.1765
$ws.Range("K134").Value = 3879.5295
$ws.Range("M134").Value = -1344.5295

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1365.2858
$ws.Range("I16").Value = 1201.75
$ws.Range("K16").Value = 1201.75
$ws.Range("M16").Value = -914.75
# Row 22
$ws.Range("H22").Value = 827.7143
$ws.Range("I22").Value = 724.25
$ws.Range("K22").Value = 724.25
$ws.Range("M22").Value = -374.25
# Row 31
$ws.Range("H31").Value = 1796.5555
$ws.Range("I31").Value = 1725.7646
$ws.Range("K31").Value = 1725.7646
$ws.Range("M31").Value = -1430.7646
# Row 34
$ws.Range("H34").Value = 1796.5555
$ws.Range("I34").Value = 1725.7646
$ws.Range("K34").Value = 1725.7646
$ws.Range("M34").Value = -1523.7646
# Row 94
$ws.Range("H94").Value = 1464.4546
$ws.Range("I94").Value = 1004.4
$ws.Range("J94").Value = 1847.8334
$ws.Range("K94").Value = 1004.4
$ws.Range("L94").Value = 1847.8334
$ws.Range("M94").Value = -553.4
$ws.Range("N94").Value = -2749.8334
# Row 113
$ws.Range("H113").Value = 1365.2858
$ws.Range("I113").Value = 1201.75
$ws.Range("K113").Value = 1201.75
$ws.Range("M113").Value = 968.25
# Row 132
$ws.Range("H132").Value = 8234.666999999999
$ws.Range("I132").Value = 7412
$ws.Range("K132").Value = 22236
$ws.Range("M132").Value = -19706

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 1536.8948
$ws.Range("J12").Value = 1592.9286
$ws.Range("L12").Value = 4778.7858
$ws.Range("N12").Value = -5124.7858

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 8448.23
$ws.Range("I113").Value = 8077.4
$ws.Range("J113").Value = 8680
$ws.Range("K113").Value = 8077.4
$ws.Range("L113").Value = 8680
$ws.Range("M113").Value = -5907.4
$ws.Range("N113").Value = -13020
# Row 132
$ws.Range("H132").Value = 5451.0728
$ws.Range("I132").Value = 4707.478
$ws.Range("K132").Value = 14122.434
$ws.Range("M132").Value = -11592.434

$ws = $wb.Worksheets.Item("LTW")
# Row 47
$ws.Range("H47").Value = 34493.5
$ws.Range("J47").Value = 34495
$ws.Range("L47").Value = 34495
$ws.Range("N47").Value = -35475
# Row 52
$ws.Range("H52").Value = 34493.5
$ws.Range("J52").Value = 34495
$ws.Range("L52").Value = 34495
$ws.Range("N52").Value = -34961
# Row 61
$ws.Range("H61").Value = 1125.75
$ws.Range("I61").Value = 501
$ws.Range("K61").Value = 501
$ws.Range("M61").Value = -299
# Row 100
$ws.Range("H100").Value = 7857.4165
$ws.Range("I100").Value = 5872.25
$ws.Range("J100").Value = 8850
$ws.Range("K100").Value = 5872.25
$ws.Range("L100").Value = 8850
$ws.Range("M100").Value = -5331.25
$ws.Range("N100").Value = -9932
# Row 113
$ws.Range("H113").Value = 1125.75
$ws.Range("I113").Value = 501
$ws.Range("K113").Value = 501
$ws.Range("M113").Value = 1669
# Row 132
$ws.Range("H132").Value = 2218.923
$ws.Range("I132").Value = 1674.1111
$ws.Range("K132").Value = 5022.3333
$ws.Range("M132").Value = -2492.3333
# Row 136
$ws.Range("H136").Value = 4208.2144
$ws.Range("I136").Value = 3789.2856
$ws.Range("K136").Value = 11367.8568
$ws.Range("M136").Value = -8817.856800000001

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2277.077
$ws.Range("I122").Value = 1538.75
$ws.Range("K122").Value = 4616.25
$ws.Range("M122").Value = -2166.25
# Row 126
$ws.Range("H126").Value = 3953.3572
$ws.Range("I126").Value = 4641.8
$ws.Range("J126").Value = 3159
$ws.Range("K126").Value = 13925.4
$ws.Range("L126").Value = 9477
$ws.Range("M126").Value = -11455.4
$ws.Range("N126").Value = -14417
# Row 136
$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1500
$ws.Range("K136").Value = 4500
$ws.Range("M136").Value = -1950
